# Redid provider indicators (renamed Provider_PPSA -> Provider_MPSA and
# recomputed its summary statistics); column width for the later sheets
# also ends up recalculated as a side-effect.

$wb = $excel.ActiveWorkbook

# Map of worksheet index -> new mean/std for the AB column ("Provider_MPSA")
# Grouped by the "count" (total observations) value found in each sheet's
# row 2 (AB2). In every year exactly 12 providers are positive for this
# indicator, so mean = 12/count and std is the corresponding sample std.
$newStats = @{
    1  = @("0.116504854368932", "0.3223982093132047")
    2  = @("0.116504854368932", "0.3223982093132047")
    3  = @("0.116504854368932", "0.3223982093132047")
    4  = @("0.116504854368932", "0.3223982093132047")
    5  = @("0.116504854368932", "0.3223982093132047")
    6  = @("0.116504854368932", "0.3223982093132047")
    7  = @("0.1176470588235294", "0.3237808098282633")
    8  = @("0.1176470588235294", "0.3237808098282633")
    9  = @("0.1188118811881188", "0.325180833164296")
    10 = @("0.1188118811881188", "0.325180833164296")
    11 = @("0.1188118811881188", "0.325180833164296")
    12 = @("0.1188118811881188", "0.325180833164296")
    13 = @("0.12", "0.3265986323710904")
    14 = @("0.12", "0.3265986323710904")
    15 = @("0.1212121212121212", "0.3280345698783139")
    16 = @("0.1212121212121212", "0.3280345698783139")
    17 = @("0.1212121212121212", "0.3280345698783139")
    18 = @("0.1263157894736842", "0.3339672956073309")
}

# Sheets whose AB column's best-fit width grows from 19 to 20 characters.
$widthSheets = @(13, 14, 15, 16, 17, 18)

for ($i = 1; $i -le 18; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Rename the "Provider_PPSA" column/header to "Provider_MPSA" everywhere.
    # Editing the header cell through the table column renames the ListObject
    # column (and therefore every table definition) consistently.
    $ws.Range("AB1").Value = "Provider_MPSA"

    # Update the describe()-style summary statistics for the renamed column.
    $stats = $newStats[$i]
    $ws.Range("AB4").Value = $stats[0]
    $ws.Range("AB5").Value = $stats[1]
    $ws.Range("AB8").Value = 0
    $ws.Range("AB9").Value = 0

    if ($widthSheets -contains $i) {
        $ws.Columns.Item(28).ColumnWidth = 19.17
    }
}
